$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1104.591
$ws.Range("I18").Value = 1135
$ws.Range("J18").Value = 1001.2
$ws.Range("K18").Value = 1135
$ws.Range("L18").Value = 1001.2
$ws.Range("M18").Value = -851
$ws.Range("N18").Value = -1569.2
# Row 86
$ws.Range("H86").Value = 168917.17
$ws.Range("I86").Value = 500751.5
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 500751.5
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -499628.5
$ws.Range("N86").Value = -5246
# Row 89
$ws.Range("H89").Value = 168917.17
$ws.Range("I89").Value = 500751.5
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 2503757.5
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -2498141.5
$ws.Range("N89").Value = -26232
# Row 98
$ws.Range("H98").Value = 33840.09
$ws.Range("I98").Value = 4605.7085
$ws.Range("J98").Value = 111798.445
$ws.Range("K98").Value = 4605.7085
$ws.Range("L98").Value = 111798.445
$ws.Range("M98").Value = -3107.7085
$ws.Range("N98").Value = -114794.445
# Row 107
$ws.Range("H107").Value = 551.4211
$ws.Range("I107").Value = 305.57144
$ws.Range("J107").Value = 1239.8
$ws.Range("K107").Value = 305.57144
$ws.Range("L107").Value = 1239.8
$ws.Range("M107").Value = 1614.42856
$ws.Range("N107").Value = -5079.8
# Row 122
$ws.Range("H122").Value = 33840.09
$ws.Range("I122").Value = 4605.7085
$ws.Range("J122").Value = 111798.445
$ws.Range("K122").Value = 13817.1255
$ws.Range("L122").Value = 335395.335
$ws.Range("M122").Value = -11367.1255
$ws.Range("N122").Value = -340295.335
# Row 129
$ws.Range("H129").Value = 24122.418
$ws.Range("I129").Value = 586.2857
$ws.Range("J129").Value = 35484.69
$ws.Range("K129").Value = 1758.8571
$ws.Range("L129").Value = 106454.07
$ws.Range("M129").Value = 3241.1429
$ws.Range("N129").Value = -116454.07
# Row 137
$ws.Range("H137").Value = 1627
$ws.Range("I137").Value = 1069.4166
$ws.Range("J137").Value = 3299.75
$ws.Range("K137").Value = 3208.2498
$ws.Range("L137").Value = 9899.25
$ws.Range("M137").Value = -658.2498000000001
$ws.Range("N137").Value = -14999.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 629785.7
$ws.Range("I88").Value = 1003625
$ws.Range("J88").Value = 6720.1665
$ws.Range("K88").Value = 1003625
$ws.Range("L88").Value = 6720.1665
$ws.Range("M88").Value = -1003219
$ws.Range("N88").Value = -7532.1665
# Row 91
$ws.Range("H91").Value = 629785.7
$ws.Range("I91").Value = 1003625
$ws.Range("J91").Value = 6720.1665
$ws.Range("K91").Value = 1003625
$ws.Range("L91").Value = 6720.1665
$ws.Range("M91").Value = -1002221
$ws.Range("N91").Value = -9528.166499999999
# Row 119
$ws.Range("H119").Value = 32142.857
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 32142.857
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 32142.857
$ws.Range("N119").Value = -41818.857

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3142.8572
$ws.Range("I86").Value = 2266.6667
$ws.Range("J86").Value = 3800
$ws.Range("K86").Value = 2266.6667
$ws.Range("L86").Value = 3800
$ws.Range("M86").Value = -1143.6667
$ws.Range("N86").Value = -6046
# Row 89
$ws.Range("H89").Value = 3142.8572
$ws.Range("I89").Value = 2266.6667
$ws.Range("J89").Value = 3800
$ws.Range("K89").Value = 11333.3335
$ws.Range("L89").Value = 19000
$ws.Range("M89").Value = -5717.333500000001
$ws.Range("N89").Value = -30232

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 74
$ws.Range("H74").Value = 20042.5
$ws.Range("I74").Value = 10285
# Row 77
$ws.Range("H77").Value = 20042.5
$ws.Range("I77").Value = 10285

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 41
$ws.Range("H41").Value = 1170
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 1170
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 3510
$ws.Range("N41").Value = -4186
# Row 69
$ws.Range("H69").Value = 1307
$ws.Range("I69").Value = 167.5
$ws.Range("J69").Value = 2066.6667
$ws.Range("K69").Value = 502.5
$ws.Range("L69").Value = 6200.000100000001
$ws.Range("M69").Value = 308.5
$ws.Range("N69").Value = -7822.000100000001
# Row 72
$ws.Range("H72").Value = 1307
$ws.Range("I72").Value = 167.5
$ws.Range("J72").Value = 2066.6667
$ws.Range("K72").Value = 1507.5
$ws.Range("L72").Value = 18600.0003
$ws.Range("M72").Value = 2548.5
$ws.Range("N72").Value = -26712.0003
# Row 93
$ws.Range("H93").Value = 2953.5715
$ws.Range("I93").Value = 350
$ws.Range("J93").Value = 3153.8462
$ws.Range("K93").Value = 1050
$ws.Range("L93").Value = 9461.5386
$ws.Range("M93").Value = 822
$ws.Range("N93").Value = -13205.5386
# Row 129
$ws.Range("H129").Value = 13890508
$ws.Range("I129").Value = 1010
$ws.Range("J129").Value = 18520342
$ws.Range("K129").Value = 3030
$ws.Range("L129").Value = 55561026
$ws.Range("M129").Value = 1970
$ws.Range("N129").Value = -55571026
# Row 131
$ws.Range("H131").Value = 766.0700000000001
$ws.Range("I131").Value = 446.22223
$ws.Range("J131").Value = 797.7033
$ws.Range("K131").Value = 1338.66669
$ws.Range("L131").Value = 2393.1099
$ws.Range("M131").Value = 3701.33331
$ws.Range("N131").Value = -12473.1099

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 103
$ws.Range("H103").Value = 19075.25
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 19075.25
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 19075.25
$ws.Range("N103").Value = -21419.25
# Row 107
$ws.Range("H107").Value = 392.6
$ws.Range("I107").Value = 470.26315
$ws.Range("J107").Value = 258.45456
$ws.Range("K107").Value = 470.26315
$ws.Range("L107").Value = 258.45456
$ws.Range("M107").Value = 1449.73685
$ws.Range("N107").Value = -4098.45456

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 18520424
$ws.Range("I61").Value = 1905.8
$ws.Range("J61").Value = 111113016
$ws.Range("K61").Value = 1905.8
$ws.Range("L61").Value = 111113016
$ws.Range("M61").Value = -1703.8
$ws.Range("N61").Value = -111113420
# Row 113
$ws.Range("H113").Value = 18520424
$ws.Range("I113").Value = 1905.8
$ws.Range("J113").Value = 111113016
$ws.Range("K113").Value = 1905.8
$ws.Range("L113").Value = 111113016
$ws.Range("M113").Value = 264.2
$ws.Range("N113").Value = -111117356
# Row 122
$ws.Range("H122").Value = 3173.4092
$ws.Range("I122").Value = 3267
$ws.Range("J122").Value = 2972.8572
$ws.Range("K122").Value = 9801
$ws.Range("L122").Value = 8918.571599999999
$ws.Range("M122").Value = -7351
$ws.Range("N122").Value = -13818.5716

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4500
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -5248
# Row 65
$ws.Range("H65").Value = 4500
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -26240
# Row 119
$ws.Range("H119").Value = 66200
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 66200
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 58100
$ws.Range("N119").Value = -75876
# Row 126
$ws.Range("H126").Value = 52638284
$ws.Range("I126").Value = 111124056
$ws.Range("J126").Value = 1090.4
$ws.Range("K126").Value = 333372168
$ws.Range("L126").Value = 3271.2
$ws.Range("M126").Value = -333369698
$ws.Range("N126").Value = -8211.200000000001
# Row 138
$ws.Range("H138").Value = 54354.832
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 54354.832
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 54354.832
$ws.Range("N138").Value = -64634.832
